# Regenerate merged AHB files
# 1) Rename the "_old"/"_new" header-suffix columns to "_FV2410"/"_FV2504"
# 2) Wrap the data range in a native Excel Table (ListObject)
# 3) Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) -----------------------------------
$ws.Range("A1").Value = "Segmentname_FV2410"
$ws.Range("B1").Value = "Segmentgruppe_FV2410"
$ws.Range("C1").Value = "Segment_FV2410"
$ws.Range("D1").Value = "Datenelement_FV2410"
$ws.Range("E1").Value = "Segment ID_FV2410"
$ws.Range("F1").Value = "Code_FV2410"
$ws.Range("G1").Value = "Qualifier_FV2410"
$ws.Range("H1").Value = "Beschreibung_FV2410"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("J1").Value = "Bedingung_FV2410"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2504"
$ws.Range("M1").Value = "Segmentgruppe_FV2504"
$ws.Range("N1").Value = "Segment_FV2504"
$ws.Range("O1").Value = "Datenelement_FV2504"
$ws.Range("P1").Value = "Segment ID_FV2504"
$ws.Range("Q1").Value = "Code_FV2504"
$ws.Range("R1").Value = "Qualifier_FV2504"
$ws.Range("S1").Value = "Beschreibung_FV2504"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2504"
$ws.Range("U1").Value = "Bedingung_FV2504"

# --- 2. Turn A1:U71 into an Excel table (adds xl/tables/table1.xml,
#        wires up sheet1's tableParts + the relationship automatically) --
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U71"), $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (pane split after row 1) ------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
